$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '42.787.30'
$ws.Range("E2").Value = '  -1.12%  '
$ws.Range("D3").Value = '2.274.92'
$ws.Range("E3").Value = '  -0.93%  '
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").Value = '''250.02'
$ws.Range("E5").Value = '  -1.01%  '
$ws.Range("D6").Value = '''0.641'
$ws.Range("E6").Value = '  -0.30%  '
$ws.Range("D7").Value = '''79.06'
$ws.Range("E7").Value = '  +7.85%  '
$ws.Range("E8").Value = '  +0.04%  '
$ws.Range("D9").Value = '''0.646'
$ws.Range("E9").Value = '  -3.32%  '
$ws.Range("D10").Value = '''41.34'
$ws.Range("E10").Value = '  +3.88%  '
$ws.Range("D11").Value = '''0.0970'
$ws.Range("E11").Value = '  -1.42%  '
$ws.Range("D12").Value = '''7.37'
$ws.Range("E12").Value = '  -2.05%  '
$ws.Range("D13").Value = '''0.106'
$ws.Range("E13").Value = '  +0.28%  '
$ws.Range("D14").Value = '2.612.12'
$ws.Range("E14").Value = '  -0.62%  '
$ws.Range("D15").Value = '''15.09'
$ws.Range("E15").Value = '  -0.70%  '
$ws.Range("D16").Value = '''0.870'
$ws.Range("E16").Value = '  -3.21%  '
$ws.Range("D17").Value = '2.279.88'
$ws.Range("E17").Value = '  -0.45%  '
$ws.Range("D18").Value = '42.660.70'
$ws.Range("E18").Value = '  -1.20%  '
$ws.Range("D19").Value = '0.0₃0996'
$ws.Range("E19").Value = '  -1.54%  '
$ws.Range("D20").Value = '''6.22'
$ws.Range("E20").Value = '  -3.15%  '
$ws.Range("E21").Value = '  -2.40%  '
$ws.Range("D22").Value = '''232.15'
$ws.Range("E22").Value = '  -2.65%  '
$ws.Range("D23").Value = '''2.17'
$ws.Range("E23").Value = '  -1.14%  '
$ws.Range("D24").Value = '''3.80'
$ws.Range("E24").Value = '  -2.80%  '
$ws.Range("E25").Value = '  +0.00%  '
$ws.Range("D26").Value = '''11.40'
$ws.Range("E26").Value = '  -4.41%  '
$ws.Range("D27").Value = '''2.34'
$ws.Range("E27").Value = '  -5.09%  '
$ws.Range("D28").Value = '''2.23'
$ws.Range("E28").Value = '  +1.31%  '
$ws.Range("D29").Value = '''169.58'
$ws.Range("E29").Value = '  +0.71%  '
$ws.Range("D30").Value = '''6.78'
$ws.Range("E30").Value = '  +7.10%  '
$ws.Range("E31").Value = '  -2.15%  '
$ws.Range("E32").Value = '  +4.89%  '
$ws.Range("D33").Value = '''0.123'
$ws.Range("E33").Value = '  -4.61%  '
$ws.Range("D34").Value = '''30.51'
$ws.Range("D35").Value = '''0.127'
$ws.Range("E35").Value = '  +0.01%  '
$ws.Range("D36").Value = '''4.58'
$ws.Range("E36").Value = '  -5.45%  '
$ws.Range("D37").Value = '''4.77'
$ws.Range("E37").Value = '  -0.88%  '
$ws.Range("E38").Value = '  -3.11%  '
$ws.Range("D39").Value = '''13.56'
$ws.Range("E39").Value = '  -0.54%  '
$ws.Range("E40").Value = '  -3.53%  '
$ws.Range("D41").Value = '''5.97'
$ws.Range("E41").Value = '  -3.06%  '
$ws.Range("D42").Value = '''115.74'
$ws.Range("E42").Value = '  +18.11%  '
$ws.Range("E43").Value = '  -2.23%  '
$ws.Range("D44").Value = '''61.53'
$ws.Range("E44").Value = '  -1.22%  '
$ws.Range("D45").Value = '''8.90'
$ws.Range("E45").Value = '  -3.66%  '
$ws.Range("E46").Value = '  -1.82%  '
$ws.Range("D47").Value = '''4.58'
$ws.Range("E47").Value = '  -7.16%  '
$ws.Range("E48").Value = '  -0.10%  '
$ws.Range("E49").Value = '  -3.80%  '
$ws.Range("D50").Value = '''1.18'
$ws.Range("E50").Value = '  -2.44%  '
$ws.Range("E51").Value = '  -2.29%  '
